# Adds USB connector (J3) to the Bonk Daddy 2HP LCSC BOM table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BonkDaddy2HP_Board")
$lo = $ws.ListObjects.Item(1)

# Insert a new row above the old row 4 (Q3 / 2N7002 / C8545), shifting the
# rest of the BOM down by one row.
$ws.Rows("4:4").Insert()

# Populate the new row with the USB connector entry.
$ws.Range("A4").Value = "J3 "
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "USB_A"
$ws.Range("D4").Value = "C718031"

# Grow the query table / autofilter to cover the new row.
$lo.Resize($ws.Range("A1:D15"))

# Keep the external-data defined name in sync with the new table extent.
$wb.Names.Item("ExternalData_1").RefersTo = "=BonkDaddy2HP_Board!`$A`$1:`$D`$15"
